# Adding Project To JENKINS
# ---------------------------------------------------------------------------
# Updates the Naukri login credentials on row 2 and appends a second data
# row (row 3) with the same login but a narrower "Location" value, mirroring
# the source worksheet's existing layout (hyperlinked username/password
# columns + plain numeric application/experience columns).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remember the existing cell formatting (xf/style) so it can be reapplied
# later -- this keeps cell styles identical to the originals instead of
# picking up whatever default formatting Range.Value / Hyperlinks.Add would
# otherwise introduce.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats -> hyperlinked-cell look (A2/B2)
$wb.Application.CutCopyMode = $false

$ws.Range("C2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # xlPasteFormats -> plain data-cell look (C2/D2/E2/F2)
$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Cell values. The order in which brand-new text is introduced matters for
# how the workbook's shared-string table is laid out, so update A2, then the
# new D3 location text, then B2, matching the source document's layout.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "siddhesh.vedre12@gmail.com"
$ws.Range("D3").Value = "Mumbai, Navi Mumbai"
$ws.Range("B2").Value = "naukri@Test123"
$ws.Range("E2").Value = 2

$ws.Range("A3").Value = "siddhesh.vedre12@gmail.com"
$ws.Range("B3").Value = "naukri@Test123"
$ws.Range("C3").Value = "Automation Test Engineer, Java, Selenium "
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 3

# ---------------------------------------------------------------------------
# Apply formatting to the new row before creating its hyperlinks, so row 3
# already shares row 2's base style and the hyperlink helper below doesn't
# need to invent a second, differently-based style variant.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$ws.Range("Z2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Hyperlinks: drop the two stale ones (they pointed at the old credentials)
# and recreate mailto links for both rows.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:siddhesh.vedre12@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "siddhesh.vedre12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:naukri@Test123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "naukri@Test123")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:siddhesh.vedre12@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "siddhesh.vedre12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:naukri@Test123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "naukri@Test123")

# Re-apply the original hyperlink-cell formatting (Hyperlinks.Add re-styles
# the cells it touches).
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Clean up the scratch cells used to stash formatting.
$ws.Range("Z1:Z2").Clear()

# Match the saved selection state from the source workbook.
[void]$ws.Range("E8").Select()
